$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh reshuffles which date/measurements land on which
# row (rows 2-14), while columns A,B,C,E,F,G,H,I,N,Q,R stay put.
# Capture the "before" state of the columns that move (D,J,K,L,M,O,P)
# for every data row, then write them back out according to the new
# row order.

$cols = @("D", "J", "K", "L", "M", "O", "P")

$before = @{}
for ($r = 2; $r -le 14; $r++) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $before[$r] = $rowData
}

# Map: new row -> row whose values it should now hold
$rowMap = @{
    2  = 11
    3  = 3
    4  = 7
    5  = 9
    6  = 13
    7  = 8
    8  = 14
    9  = 6
    10 = 2
    11 = 10
    12 = 4
    13 = 5
    14 = 12
}

foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $src = $before[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $src[$col]
    }
}
